$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update absorption voltage values from 55.6 to 55.2 VDC (rows 3-19 and 24-26)
$ws.Range("G3:G19").Value = 55.2
$ws.Range("G24:G26").Value = 55.2

# Update the active cell selection to H25
$ws.Range("H25").Select()
